$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.638.30'
$ws.Range("E2").Value = '  -0.80%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.616.57'
$ws.Range("E3").Value = '  -0.70%  '
$ws.Range("E4").Value = '  -0.57%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '209.23'
$ws.Range("E5").Value = '  -0.98%  '
$ws.Range("E6").Value = '  -1.18%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.989'
$ws.Range("E7").Value = '  -0.53%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.05'
$ws.Range("E8").Value = '  -0.84%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.254'
$ws.Range("E9").Value = '  -1.31%  '
$ws.Range("E10").Value = '  -1.26%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0873'
$ws.Range("E11").Value = '  -0.81%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.847.21'
$ws.Range("E12").Value = '  -0.64%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.624.13'
$ws.Range("E13").Value = '  -0.26%  '
$ws.Range("E14").Value = '  -1.52%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.555'
$ws.Range("E15").Value = '  -1.40%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.50'
$ws.Range("E16").Value = '  -1.03%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.658.43'
$ws.Range("E17").Value = '  -0.71%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '227.10'
$ws.Range("E18").Value = '  -1.40%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.63'
$ws.Range("E19").Value = '  +1.50%  '
$ws.Range("E20").Value = '  -0.94%  '
$ws.Range("E21").Value = '  -0.48%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.30'
$ws.Range("E22").Value = '  -1.29%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.03'
$ws.Range("E23").Value = '  -2.76%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.04'
$ws.Range("E24").Value = '  -0.62%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.36'
$ws.Range("E25").Value = '  +0.10%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.87'
$ws.Range("E26").Value = '  -1.05%  '
$ws.Range("E27").Value = '  -0.55%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.39'
$ws.Range("E28").Value = '  -1.31%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.991'
$ws.Range("E29").Value = '  -0.51%  '
$ws.Range("E30").Value = '  -0.78%  '
$ws.Range("E31").Value = '  -0.61%  '
$ws.Range("E32").Value = '  -1.23%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.06'
$ws.Range("E33").Value = '  -0.40%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.388.54'
$ws.Range("E34").Value = '  -1.11%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.58'
$ws.Range("E35").Value = '  +1.37%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.997'
$ws.Range("E36").Value = '  -1.28%  '
$ws.Range("E37").Value = '  -1.35%  '
$ws.Range("E38").Value = '  +0.38%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.554'
$ws.Range("E39").Value = '  -0.88%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.840'
$ws.Range("E40").Value = '  -3.02%  '
$ws.Range("E41").Value = '  -1.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.989'
$ws.Range("E43").Value = '  -0.49%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '65.42'
$ws.Range("E44").Value = '  -1.73%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.35'
$ws.Range("E45").Value = '  -2.87%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.756.67'
$ws.Range("E46").Value = '  -0.76%  '
$ws.Range("E47").Value = '  -3.45%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.51'
$ws.Range("E48").Value = '  -0.26%  '
$ws.Range("E49").Value = '  +1.10%  '
$ws.Range("E50").Value = '  -0.66%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.52'
$ws.Range("E51").Value = '  +0.87%  '
